$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before column H (currently: A..J = id, name, owner, quantity,
# face_value, currency, total, date, legislator_name, legislator_id).
# The new column becomes "property_category" with constant value "stock",
# and the former H/I/J (date/legislator_name/legislator_id) shift right to I/J/K.
$ws.Columns("H:H").Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H4").Value = "stock"

# Fix the quantity value in D2: remove the full-width comma from "26，773" -> "26773"
# Keep it stored as text (not a number) to match the original cell's string type.
# (A direct .Value assignment of a numeric-looking string gets auto-converted to a
# number, or needs a quote-prefix that introduces a new cell style; going through a
# text-producing formula and then pasting as a value avoids both problems.)
$ws.Range("D2").Formula = '="26773"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
